$wb = $excel.ActiveWorkbook

# --- TagSave sheet: append two more rows (105, 106) ---
$tagSave = $wb.Worksheets.Item("TagSave")
$tagSave.Range("A105").Value = "30VQmLjkKU"
$tagSave.Range("A106").Value = "Fs18cQJDnU"

# --- CitySave sheet: append three more rows (2, 3, 4) ---
$citySave = $wb.Worksheets.Item("CitySave")
$citySave.Range("A2").Value = "KE"
$citySave.Range("A3").Value = "Pt"
$citySave.Range("A4").Value = "jX"
